$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (date) values between row 2 and row 4
$ws.Range("D2").Value = 44253
$ws.Range("D4").Value = 44250

# Swap the Volumen values between row 2 and row 4
$ws.Range("M2").Value = 160
$ws.Range("M4").Value = 200
